# Update the loading_percent results sheet for Case_3_10 (380 kV case).
# Columns updated (0-indexed data columns from header row 1): B, C, D, F, G, H, I, J, L, N
# Columns E, K, M stay at 0 and are left untouched, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows correspond to worksheet rows 2..25 (the "case index" column A holds 0..23).
$columns = @("B", "C", "D", "F", "G", "H", "I", "J", "L", "N")

$data = @(
    @(19.72533898087631, 7.955519074909731, 7.778649110584195, 37.63345625936837, 43.34716879409017, 17.78302384349977, 26.96349738015721, 9.965697282616896, 12.25303487422587, 18.57794907394481),
    @(19.30689682440947, 7.505145175652517, 7.779906183598751, 37.65022079033169, 43.2698189046213, 17.82296733596651, 27.05312787994748, 9.987062619886832, 12.24283065015564, 18.64248782568058),
    @(19.0491330102792, 7.212995233458619, 7.781190844707325, 37.67039197622758, 43.236648158842, 17.85114529058134, 27.11414039461181, 10.00090491940092, 12.23836574768118, 18.68400535037199),
    @(18.94403188811371, 7.090055297196023, 7.781843804219828, 37.68109167760533, 43.22673628102135, 17.86354447854191, 27.14050260347101, 10.00672829983117, 12.23700097048673, 18.70140091660891),
    @(18.92658076470734, 7.069407885212675, 7.78196006046607, 37.68301798161709, 43.22530821786425, 17.86565864366773, 27.14497044142416, 10.00770630690513, 12.23680186614813, 18.70431828012015),
    @(19.04771562196267, 7.211352896908154, 7.781199125948911, 37.67052624207813, 43.23649988329816, 17.85130880252315, 27.11448986062177, 10.00098271581867, 12.23834549831493, 18.68423802043691),
    @(19.58132322907307, 7.803493261365737, 7.778976416195617, 37.63718453098425, 43.31752957449301, 17.79603712113672, 26.99315843908315, 9.972914108538061, 12.24914398663246, 18.59981055127801),
    @(20.61469025735193, 8.839430358935234, 7.778664900336413, 37.65029882593037, 43.58970617758781, 17.71671513177136, 26.80286010331581, 9.923592569743638, 12.2845130318366, 18.44918120081894),
    @(21.35781807825951, 9.522944594747386, 7.780874375396999, 37.70786842652545, 43.85796941228281, 17.67626488841594, 26.69234189685452, 9.890811648227984, 12.31901523334594, 18.34752203676584),
    @(21.69080567732661, 9.816902804817921, 7.782402623118006, 37.74445253211586, 43.99458885103197, 17.66175403652001, 26.64847763518751, 9.876642321440841, 12.33652705346297, 18.30321053158468),
    @(21.81604909106743, 9.925769215698461, 7.783055985764771, 37.75979684156582, 44.04839279970728, 17.65681986642091, 26.63279337836119, 9.871383080603797, 12.34341632006646, 18.28670750133054),
    @(21.7891154422922, 9.902431906782764, 7.7829119617966, 37.75642595371589, 44.03671362246308, 17.65785756783764, 26.63613000469502, 9.872511028237511, 12.34192117707732, 18.29024943498577),
    @(21.7011272770724, 9.825908429355259, 7.782454882957827, 37.74568507126284, 43.99897404135727, 17.6613368534481, 26.64716869578069, 9.876207511162386, 12.33708868815997, 18.30184727773599),
    @(21.64711753586646, 9.778716452964579, 7.78218461296718, 37.73929995357052, 43.97612601579524, 17.66354108317693, 26.6540509542386, 9.878485553684278, 12.33416213624414, 18.30898729739258),
    @(21.33594389994913, 9.503391365899356, 7.780784978643163, 37.70568643393633, 43.84933218217569, 17.67729160864175, 26.69533790606831, 9.891752557161501, 12.31790707182996, 18.35045669872063),
    @(21.14366131589797, 9.330133531084389, 7.780059930770736, 37.68772570522756, 43.77526509191051, 17.68672451425595, 26.72231123371126, 9.900081381999538, 12.30839810887612, 18.37639118667408),
    @(21.03259559546266, 9.228882592802611, 7.779692178993301, 37.67837398414892, 43.73403842990216, 17.6925161628198, 26.73842876749206, 9.904941851601288, 12.30310000315606, 18.39149010902892),
    @(20.99491384441475, 9.194326677551826, 7.779576145278587, 37.67537586441023, 43.72031670205201, 17.69453994697633, 26.74398935255559, 9.906599550816402, 12.30133565614538, 18.39663365881898),
    @(21.16417967408447, 9.348742539998479, 7.780132017577223, 37.68953637660539, 43.78300756814046, 17.6856824657286, 26.71937741748834, 9.899187528098201, 12.30939265907769, 18.3736115798688),
    @(21.72699559406608, 9.848451747173662, 7.782587116842926, 37.74879951258222, 44.01000315982399, 17.66029967368507, 26.64390119479267, 9.87511888107613, 12.33850113481841, 18.29843320883557),
    @(22.08981210771093, 10.16076342677964, 7.784626495630867, 37.79621729301865, 44.17040254921115, 17.64697967324157, 26.59997352060708, 9.860008462565483, 12.35902687900328, 18.25091246515231),
    @(21.8966677639654, 9.995384979501704, 7.783498448825513, 37.77011655710798, 44.08370262021442, 17.6537892803994, 26.62292303539375, 9.86801660978378, 12.34793564724776, 18.27612804515405),
    @(21.15490493420191, 9.340334523187249, 7.780099274204448, 37.68871473729927, 43.77950297016624, 17.68615242791475, 26.7207018952815, 9.899591414770237, 12.30894249712751, 18.37486765174109),
    @(20.33739856330205, 8.572742549254924, 7.778318828301871, 37.63833300738526, 43.50400996744541, 17.73505029504486, 26.8492145432915, 9.936326215007249, 12.2734384359815, 18.4883416863521)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $columns.Count; $j++) {
        $cellRef = $columns[$j] + $rowNum
        $ws.Range($cellRef).Value = $rowValues[$j]
    }
}

Write-Host "Updated loading_percent values for rows 2-25."
